$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: stage copies of source rows 14-18 and 32-36 to temporary rows (offset +200).
# The staging rows are cleared first so that pasting (which only transfers cells that
# actually contain data in the source) cannot leave stale leftover values behind.
# This lets the permutation/swap among rows 14-18 and 32-36 be performed safely without
# a destination row retaining any of its own previous data.
$ws.Range("A214:AY214").ClearContents() | Out-Null
$ws.Range("A14:AY14").Copy() | Out-Null
$ws.Range("A214:AY214").PasteSpecial(-4163) | Out-Null
$ws.Range("A215:AY215").ClearContents() | Out-Null
$ws.Range("A15:AY15").Copy() | Out-Null
$ws.Range("A215:AY215").PasteSpecial(-4163) | Out-Null
$ws.Range("A216:AY216").ClearContents() | Out-Null
$ws.Range("A16:AY16").Copy() | Out-Null
$ws.Range("A216:AY216").PasteSpecial(-4163) | Out-Null
$ws.Range("A217:AY217").ClearContents() | Out-Null
$ws.Range("A17:AY17").Copy() | Out-Null
$ws.Range("A217:AY217").PasteSpecial(-4163) | Out-Null
$ws.Range("A218:AY218").ClearContents() | Out-Null
$ws.Range("A18:AY18").Copy() | Out-Null
$ws.Range("A218:AY218").PasteSpecial(-4163) | Out-Null
$ws.Range("A232:AY232").ClearContents() | Out-Null
$ws.Range("A32:AY32").Copy() | Out-Null
$ws.Range("A232:AY232").PasteSpecial(-4163) | Out-Null
$ws.Range("A233:AY233").ClearContents() | Out-Null
$ws.Range("A33:AY33").Copy() | Out-Null
$ws.Range("A233:AY233").PasteSpecial(-4163) | Out-Null
$ws.Range("A234:AY234").ClearContents() | Out-Null
$ws.Range("A34:AY34").Copy() | Out-Null
$ws.Range("A234:AY234").PasteSpecial(-4163) | Out-Null
$ws.Range("A235:AY235").ClearContents() | Out-Null
$ws.Range("A35:AY35").Copy() | Out-Null
$ws.Range("A235:AY235").PasteSpecial(-4163) | Out-Null
$ws.Range("A236:AY236").ClearContents() | Out-Null
$ws.Range("A36:AY36").Copy() | Out-Null
$ws.Range("A236:AY236").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

# Step 2: copy staged rows into their final destination rows per the row permutation,
# clearing each destination first for the same reason as above.
$ws.Range("A14:AY14").ClearContents() | Out-Null
$ws.Range("A214:AY214").Copy() | Out-Null
$ws.Range("A14:AY14").PasteSpecial(-4163) | Out-Null
$ws.Range("A15:AY15").ClearContents() | Out-Null
$ws.Range("A216:AY216").Copy() | Out-Null
$ws.Range("A15:AY15").PasteSpecial(-4163) | Out-Null
$ws.Range("A16:AY16").ClearContents() | Out-Null
$ws.Range("A215:AY215").Copy() | Out-Null
$ws.Range("A16:AY16").PasteSpecial(-4163) | Out-Null
$ws.Range("A17:AY17").ClearContents() | Out-Null
$ws.Range("A218:AY218").Copy() | Out-Null
$ws.Range("A17:AY17").PasteSpecial(-4163) | Out-Null
$ws.Range("A18:AY18").ClearContents() | Out-Null
$ws.Range("A217:AY217").Copy() | Out-Null
$ws.Range("A18:AY18").PasteSpecial(-4163) | Out-Null
$ws.Range("A32:AY32").ClearContents() | Out-Null
$ws.Range("A235:AY235").Copy() | Out-Null
$ws.Range("A32:AY32").PasteSpecial(-4163) | Out-Null
$ws.Range("A33:AY33").ClearContents() | Out-Null
$ws.Range("A234:AY234").Copy() | Out-Null
$ws.Range("A33:AY33").PasteSpecial(-4163) | Out-Null
$ws.Range("A34:AY34").ClearContents() | Out-Null
$ws.Range("A233:AY233").Copy() | Out-Null
$ws.Range("A34:AY34").PasteSpecial(-4163) | Out-Null
$ws.Range("A35:AY35").ClearContents() | Out-Null
$ws.Range("A232:AY232").Copy() | Out-Null
$ws.Range("A35:AY35").PasteSpecial(-4163) | Out-Null
$ws.Range("A36:AY36").ClearContents() | Out-Null
$ws.Range("A236:AY236").Copy() | Out-Null
$ws.Range("A36:AY36").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

# Step 3: apply the Taxonsorteringsordning (column B) value updates
$ws.Range("B14").Value = 56430
$ws.Range("B15").Value = 78699
$ws.Range("B16").Value = 56430
$ws.Range("B17").Value = 56430
$ws.Range("B18").Value = 90099
$ws.Range("B32").Value = 78699
$ws.Range("B33").Value = 89539
$ws.Range("B34").Value = 90099
$ws.Range("B35").Value = 89539
$ws.Range("B36").Value = 56430

# Step 4: clear the temporary staging rows
$ws.Range("A214:AY214").ClearContents() | Out-Null
$ws.Range("A215:AY215").ClearContents() | Out-Null
$ws.Range("A216:AY216").ClearContents() | Out-Null
$ws.Range("A217:AY217").ClearContents() | Out-Null
$ws.Range("A218:AY218").ClearContents() | Out-Null
$ws.Range("A232:AY232").ClearContents() | Out-Null
$ws.Range("A233:AY233").ClearContents() | Out-Null
$ws.Range("A234:AY234").ClearContents() | Out-Null
$ws.Range("A235:AY235").ClearContents() | Out-Null
$ws.Range("A236:AY236").ClearContents() | Out-Null
